$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark these tools as "Featured" (column S) by setting value to "Y"
$ws.Range("S5:S9").Value = "Y"

# Update the visible selection/scroll position left by the editor
$excel.ActiveWindow.ScrollColumn = 14
$ws.Range("S16").Select()
